$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$cell = $ws.Range('D2')
$cell.NumberFormat = '@'
$cell.Value = '26.873.07'
$cell.Style = 'Normal'
$ws.Range('E2').Value = '  -0.64%  '

# Row 3
$cell = $ws.Range('D3')
$cell.NumberFormat = '@'
$cell.Value = '1.870.92'
$cell.Style = 'Normal'
$ws.Range('E3').Value = '  +0.28%  '

# Row 4
$ws.Range('E4').Value = '  +0.03%  '

# Row 5
$cell = $ws.Range('D5')
$cell.NumberFormat = '@'
$cell.Value = '304.85'
$cell.Style = 'Normal'
$ws.Range('E5').Value = '  -0.35%  '

# Row 6
$cell = $ws.Range('D6')
$cell.NumberFormat = '@'
$cell.Value = '0.9999'
$cell.Style = 'Normal'
$ws.Range('E6').Value = '  +0.03%  '

# Row 7
$cell = $ws.Range('D7')
$cell.NumberFormat = '@'
$cell.Value = '0.5086'
$cell.Style = 'Normal'
$ws.Range('E7').Value = '  -1.44%  '

# Row 8
$cell = $ws.Range('D8')
$cell.NumberFormat = '@'
$cell.Value = '0.3661'
$cell.Style = 'Normal'

# Row 9
$cell = $ws.Range('D9')
$cell.NumberFormat = '@'
$cell.Value = '0.07188'
$cell.Style = 'Normal'
$ws.Range('E9').Value = '  +0.39%  '

# Row 10
$cell = $ws.Range('D10')
$cell.NumberFormat = '@'
$cell.Value = '0.8928'
$cell.Style = 'Normal'
$ws.Range('E10').Value = '  +0.44%  '

# Row 11
$cell = $ws.Range('D11')
$cell.NumberFormat = '@'
$cell.Value = '20.67'
$cell.Style = 'Normal'
$ws.Range('E11').Value = '  +0.02%  '

# Row 12
$cell = $ws.Range('D12')
$cell.NumberFormat = '@'
$cell.Value = '0.07522'
$cell.Style = 'Normal'
$ws.Range('E12').Value = '  -1.12%  '

# Row 13
$cell = $ws.Range('D13')
$cell.NumberFormat = '@'
$cell.Value = '1.879.55'
$cell.Style = 'Normal'
$ws.Range('E13').Value = '  +0.52%  '

# Row 14
$cell = $ws.Range('D14')
$cell.NumberFormat = '@'
$cell.Value = '94.95'
$cell.Style = 'Normal'
$ws.Range('E14').Value = '  +5.98%  '

# Row 15
$cell = $ws.Range('D15')
$cell.NumberFormat = '@'
$cell.Value = '5.224'
$cell.Style = 'Normal'
$ws.Range('E15').Value = '  -1.47%  '

# Row 16
$ws.Range('E16').Value = '  +0.03%  '

# Row 17
$cell = $ws.Range('D17')
$cell.NumberFormat = '@'
$cell.Value = '0.000008501'
$cell.Style = 'Normal'
$ws.Range('E17').Value = '  +0.38%  '

# Row 18
$cell = $ws.Range('D18')
$cell.NumberFormat = '@'
$cell.Value = '14.19'
$cell.Style = 'Normal'
$ws.Range('E18').Value = '  +0.85%  '

# Row 19
$cell = $ws.Range('D19')
$cell.NumberFormat = '@'
$cell.Value = '0.9997'
$cell.Style = 'Normal'
$ws.Range('E19').Value = '  -0.04%  '

# Row 20
$cell = $ws.Range('D20')
$cell.NumberFormat = '@'
$cell.Value = '26.932.15'
$cell.Style = 'Normal'
$ws.Range('E20').Value = '  -0.51%  '

# Row 21
$cell = $ws.Range('D21')
$cell.NumberFormat = '@'
$cell.Value = '5.014'
$cell.Style = 'Normal'
$ws.Range('E21').Value = '  -0.30%  '

# Row 22
$cell = $ws.Range('D22')
$cell.NumberFormat = '@'
$cell.Value = '2.114.44'
$cell.Style = 'Normal'
$ws.Range('E22').Value = '  -0.82%  '

# Row 23
$ws.Range('E23').Value = '  -1.52%  '

# Row 24
$cell = $ws.Range('D24')
$cell.NumberFormat = '@'
$cell.Value = '6.391'
$cell.Style = 'Normal'
$ws.Range('E24').Value = '  -1.13%  '

# Row 25
$cell = $ws.Range('D25')
$cell.NumberFormat = '@'
$cell.Value = '148.26'
$cell.Style = 'Normal'
$ws.Range('E25').Value = '  +0.43%  '

# Row 26
$cell = $ws.Range('D26')
$cell.NumberFormat = '@'
$cell.Value = '1.779'
$cell.Style = 'Normal'
$ws.Range('E26').Value = '  -3.19%  '

# Row 27
$cell = $ws.Range('D27')
$cell.NumberFormat = '@'
$cell.Value = '17.87'
$cell.Style = 'Normal'
$ws.Range('E27').Value = '  -0.56%  '

# Row 28
$cell = $ws.Range('D28')
$cell.NumberFormat = '@'
$cell.Value = '2.090'
$cell.Style = 'Normal'
$ws.Range('E28').Value = '  -0.47%  '

# Row 29
$cell = $ws.Range('D29')
$cell.NumberFormat = '@'
$cell.Value = '113.28'
$cell.Style = 'Normal'
$ws.Range('E29').Value = '  +0.44%  '

# Row 30
$cell = $ws.Range('D30')
$cell.NumberFormat = '@'
$cell.Value = '4.700'
$cell.Style = 'Normal'
$ws.Range('E30').Value = '  +0.85%  '

# Row 31
$cell = $ws.Range('D31')
$cell.NumberFormat = '@'
$cell.Value = '4.732'
$cell.Style = 'Normal'
$ws.Range('E31').Value = '  +1.08%  '

# Row 32
$cell = $ws.Range('D32')
$cell.NumberFormat = '@'
$cell.Value = '0.09138'
$cell.Style = 'Normal'
$ws.Range('E32').Value = '  -0.08%  '

# Row 33
$cell = $ws.Range('D33')
$cell.NumberFormat = '@'
$cell.Value = '0.05077'
$cell.Style = 'Normal'
$ws.Range('E33').Value = '  -0.90%  '

# Row 34
$cell = $ws.Range('D34')
$cell.NumberFormat = '@'
$cell.Value = '0.7495'
$cell.Style = 'Normal'
$ws.Range('E34').Value = '  +3.24%  '

# Row 35
$cell = $ws.Range('D35')
$cell.NumberFormat = '@'
$cell.Value = '2.980'
$cell.Style = 'Normal'
$ws.Range('E35').Value = '  -2.89%  '

# Row 36
$cell = $ws.Range('D36')
$cell.NumberFormat = '@'
$cell.Value = '1.156'
$cell.Style = 'Normal'
$ws.Range('E36').Value = '  -0.03%  '

# Row 37
$cell = $ws.Range('D37')
$cell.NumberFormat = '@'
$cell.Value = '3.234'
$cell.Style = 'Normal'
$ws.Range('E37').Value = '  +5.61%  '

# Row 38
$cell = $ws.Range('D38')
$cell.NumberFormat = '@'
$cell.Value = '2.528'
$cell.Style = 'Normal'
$ws.Range('E38').Value = '  +1.23%  '

# Row 39
$ws.Range('B39').Value = 'VeChain'
$ws.Range('C39').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$cell = $ws.Range('D39')
$cell.NumberFormat = '@'
$cell.Value = '0.01995'
$cell.Style = 'Normal'
$ws.Range('E39').Value = '  -1.97%  '

# Row 40
$ws.Range('B40').Value = 'TheSandbox'
$ws.Range('C40').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$cell = $ws.Range('D40')
$cell.NumberFormat = '@'
$cell.Value = '0.5585'
$cell.Style = 'Normal'
$ws.Range('E40').Value = '  +4.78%  '

# Row 41
$cell = $ws.Range('D41')
$cell.NumberFormat = '@'
$cell.Value = '1.073'
$cell.Style = 'Normal'
$ws.Range('E41').Value = '  -0.04%  '

# Row 42
$cell = $ws.Range('D42')
$cell.NumberFormat = '@'
$cell.Value = '6.625'
$cell.Style = 'Normal'
$ws.Range('E42').Value = '  +1.52%  '

# Row 43
$cell = $ws.Range('D43')
$cell.NumberFormat = '@'
$cell.Value = '115.94'
$cell.Style = 'Normal'
$ws.Range('E43').Value = '  -0.43%  '

# Row 44
$cell = $ws.Range('D44')
$cell.NumberFormat = '@'
$cell.Value = '8.586'
$cell.Style = 'Normal'
$ws.Range('E44').Value = '  +3.55%  '

# Row 45
$cell = $ws.Range('D45')
$cell.NumberFormat = '@'
$cell.Value = '0.1476'
$cell.Style = 'Normal'
$ws.Range('E45').Value = '  +0.69%  '

# Row 46
$cell = $ws.Range('D46')
$cell.NumberFormat = '@'
$cell.Value = '0.4760'
$cell.Style = 'Normal'
$ws.Range('E46').Value = '  +2.68%  '

# Row 47
$cell = $ws.Range('D47')
$cell.NumberFormat = '@'
$cell.Value = '0.9993'
$cell.Style = 'Normal'
$ws.Range('E47').Value = '  +0.00%  '

# Row 48
$ws.Range('E48').Value = '  +1.18%  '

# Row 49
$cell = $ws.Range('D49')
$cell.NumberFormat = '@'
$cell.Value = '1.567'
$cell.Style = 'Normal'
$ws.Range('E49').Value = '  -0.24%  '

# Row 50
$cell = $ws.Range('D50')
$cell.NumberFormat = '@'
$cell.Value = '36.96'
$cell.Style = 'Normal'
$ws.Range('E50').Value = '  +1.16%  '

# Row 51
$cell = $ws.Range('D51')
$cell.NumberFormat = '@'
$cell.Value = '63.12'
$cell.Style = 'Normal'
$ws.Range('E51').Value = '  -0.74%  '
